$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "ActivePowerSummary"
$ws.Range("B12").Value = "Мощность"
$ws.Range("A13").Select()
